# Generate Report for Handoff
# --------------------------------------------------------------------------
# The report (3 sheets: Overview, zh-cn, de-de) tracked two source files
# (c5a84f31... and dcf3d87c...). A new handoff run replaces them with two
# different source files (fa48051e... and ffffc34c901e...), flips the status
# from "Handed back: in sync with en-US" to "Ready for handoff", refreshes
# the handoff timestamps, and - because both locales now share a single
# handoff package - collapses the separate "Latest Target File" / "Latest
# Handback File" columns (F/G) into the already-present "Latest Handoff
# File" column (D), so F/G are cleared out on the zh-cn/de-de sheets.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$md1 = "fa48051e-49a5-4ba8-b3c8-aa60303d4929.md"
$md2 = "ffffc34c901e-956a-41c6-abfc-f185a7ce177e.md"
$status = "Ready for handoff"
$dateOverview = "2016-03-24 15:27:00"
$xlfZhCn = "fa48051e-49a5-4ba8-b3c8-aa60303d4929.3cca7003702a782f4eb0e95c6c086d5ec4b260d9.zh-cn.xlf"
$xlfDeDe = "fa48051e-49a5-4ba8-b3c8-aa60303d4929.3cca7003702a782f4eb0e95c6c086d5ec4b260d9.de-de.xlf"
$eDateZhCn = "2016-03-24 15:26:55"
$hDate = "0001-01-01 00:00:00"

$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/8f1642f2db6ba224d37a6bd1893c4c015f527e70/e2e/$md1"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/8f1642f2db6ba224d37a6bd1893c4c015f527e70/e2e/$md2"
$xlfZhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30030daf38935a5f453c2469b765a666d5cd26f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhCn"
$xlfDeDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbacc48a0727acbb0e37ac04e9f43f984edc5e67/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeDe"

# --------------------------------------------------------------------------
# Overview sheet
# --------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = $status
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $dateOverview

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = $status
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $dateOverview

$ws.Hyperlinks.Add($ws.Range("A2"), $md1Url, "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("A3"), $md2Url, "", "", $md2)

# --------------------------------------------------------------------------
# zh-cn sheet
# --------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $md1
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $xlfZhCn
$ws.Range("E2").Value = $eDateZhCn
$ws.Range("F2:G2").ClearContents()
$ws.Range("H2").Value = $hDate

$ws.Range("A3").Value = $md2
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $xlfZhCn
$ws.Range("E3").Value = $eDateZhCn
$ws.Range("F3:G3").ClearContents()
$ws.Range("H3").Value = $hDate

$ws.Hyperlinks.Add($ws.Range("A2"), $md1Url, "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfZhCnUrl, "", "", $xlfZhCn)
$ws.Hyperlinks.Add($ws.Range("A3"), $md2Url, "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfZhCnUrl, "", "", $xlfZhCn)

# --------------------------------------------------------------------------
# de-de sheet
# --------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $md1
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $xlfDeDe
$ws.Range("E2").Value = $dateOverview
$ws.Range("F2:G2").ClearContents()
$ws.Range("H2").Value = $hDate

$ws.Range("A3").Value = $md2
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $xlfDeDe
$ws.Range("E3").Value = $dateOverview
$ws.Range("F3:G3").ClearContents()
$ws.Range("H3").Value = $hDate

$ws.Hyperlinks.Add($ws.Range("A2"), $md1Url, "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfDeDeUrl, "", "", $xlfDeDe)
$ws.Hyperlinks.Add($ws.Range("A3"), $md2Url, "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfDeDeUrl, "", "", $xlfDeDe)
